# Chiffres COVID-19 Valais - data corrections / additions
# Only the raw input columns (C, E, F, G, L, M) are set here; columns
# B, H, J, K are live formulas (cumulative totals) and Excel recomputes
# them automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 238 / 245: one case moved from 2020-10-16 (row238) to 2020-10-23 (row245)
$ws.Range("C238").Value = 558
$ws.Range("C245").Value = 882

# --- Row 616: updated daily case count
$ws.Range("C616").Value = 76

# --- Rows 635-646: updated daily figures
$ws.Range("C635").Value = 101
$ws.Range("M635").Value = 1

$ws.Range("C636").Value = 303

$ws.Range("C637").Value = 247

$ws.Range("E640").Value = 10

$ws.Range("C641").Value = 152
$ws.Range("E641").Value = 10

$ws.Range("C642").Value = 145
$ws.Range("F642").Value = 5
$ws.Range("G642").Value = 57

$ws.Range("C643").Value = 420
$ws.Range("F643").Value = 2
$ws.Range("G643").Value = 54

$ws.Range("C644").Value = 325
$ws.Range("F644").Value = 5
$ws.Range("G644").Value = 63

$ws.Range("C645").Value = 348
$ws.Range("E645").Value = 9
$ws.Range("F645").Value = 5

$ws.Range("C646").Value = 343
$ws.Range("E646").Value = 9
$ws.Range("F646").Value = 7
$ws.Range("G646").Value = 61

# --- Rows 647-650: previously-empty rows now filled in with new data
$ws.Range("C647").Value = 353
$ws.Range("E647").Value = 9
$ws.Range("F647").Value = 8
$ws.Range("G647").Value = 64
$ws.Range("L647").Value = 0
$ws.Range("M647").Value = 0

$ws.Range("C648").Value = 125
$ws.Range("E648").Value = 8
$ws.Range("F648").Value = 7
$ws.Range("G648").Value = 70
$ws.Range("L648").Value = 0
$ws.Range("M648").Value = 0

$ws.Range("C649").Value = 78
$ws.Range("E649").Value = 9
$ws.Range("F649").Value = 6
$ws.Range("G649").Value = 68
$ws.Range("L649").Value = 0
$ws.Range("M649").Value = 0

$ws.Range("C650").Value = 15
$ws.Range("E650").Value = 9
$ws.Range("F650").Value = 7
$ws.Range("G650").Value = 64
$ws.Range("L650").Value = 0
$ws.Range("M650").Value = 0
